# Rule: Mead Exclusive Flood Control Space - edit
# For every trace sheet in this workbook, Trace1 (column B, rows 2:37) is
# being replaced with the recomputed values that already live in Trace5
# (column F, rows 2:37) - the rule previously driving Mead outflow while
# Mead was nearly empty now keys off Flood Control Flag = 1 (not NaN),
# which changes the Trace1 series. Also update the sheet's saved selection
# to reflect the edited range (B2:B37), matching what was active when the
# edit was made in Excel.

$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    $srcValues = $ws.Range("F2:F37").Value()
    $ws.Range("B2:B37").Value = $srcValues

    $null = $ws.Activate()
    $ws.Range("B2:B37").Select() | Out-Null
}
